# display table bug fix
# Adds 13 new rows (957-969) to the "rush" regex lookup table, continuing the
# stend/stbegin pattern-pair list that already fills the sheet, and nudges a
# few cosmetic view/sizing properties to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rush")

# ---------------------------------------------------------------------------
# New table rows: column A = regex text, column B = 2 (stend) / 0 (stbegin),
# column C = "stend" / "stbegin" marker (shared with the rest of the sheet).
# ---------------------------------------------------------------------------
$rows = @(
    @("#\a)\w+PULM:",              2, "stend"),
    @("#\n(P)ULM",                 0, "stbegin"),
    @("#\a)\w+IMPRESSION:",        2, "stend"),
    @("#\n(I)MPRESSION:",          0, "stbegin"),
    @("#\a)\w+IMPRESSION:",        2, "stend"),
    @("#\n(I)MPRESSION:",          0, "stbegin"),
    @("#\a)\w+REASON:",            2, "stend"),
    @("#\n(R)EASON:",              0, "stbegin"),
    @("#\a)\w+CHEST:",             2, "stend"),
    @("#\n(C)HEST:",               0, "stbegin"),
    @("\a)\w+\C\C\C\C+:",          2, "stend"),
    @("\n\w+(\C)\C\C\C+:",         0, "stbegin"),
    @("\n(\C)\C\C\C+:",            0, "stbegin")
)

$startRow = 957
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $triple = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $triple[0]
    $ws.Cells.Item($r, 2).Value = $triple[1]
    $ws.Cells.Item($r, 3).Value = $triple[2]

    $ws.Rows.Item($r).RowHeight = 12.8
}

# First new row opens a fresh wrapped block, mirroring row 955's style.
$ws.Cells.Item($startRow, 1).WrapText = $true

# ---------------------------------------------------------------------------
# View bookkeeping: keep the window scrolled near the bottom of the table and
# move the active selection to the new first empty row beneath the data.
# ---------------------------------------------------------------------------
$lastRow = $startRow + $rows.Count - 1
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 936
$ws.Cells.Item($lastRow + 1, 1).Select()

# ---------------------------------------------------------------------------
# Column widths shrank slightly (content-driven autosize in the source file).
# Reproduce as closely as the host's width model allows.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 84.66666666666667
$ws.Columns.Item(2).ColumnWidth = 5.166666666666667
$ws.Columns.Item(3).ColumnWidth = 6.0
